$wb = $excel.ActiveWorkbook

# 1. Rename "Root" sheet to "Main root" (tests fuzzy worksheet-name matching)
$rootWs = $wb.Worksheets.Item("Root")
$rootWs.Name = "Main root"

# 2. For every sheet's visible _FilterDatabase defined name, add a sibling
#    "_xlnm._FilterDatabase_0" defined name with the very same reference.
foreach ($ws in $wb.Worksheets) {
    foreach ($n in $wb.Names) {
        if ($n.Name -eq ($ws.Name + "!_FilterDatabase") -and $n.Visible) {
            $ws.Names.Add("_xlnm._FilterDatabase_0", $n.RefersTo)
        }
    }
}

# 3. Move the active tab/window selection from "One to many rows" back to the
#    first sheet ("Main root").
$rootWs.Activate()
